$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Predictor column (C) values: wrap per-capita / similar predictors with ln(...) and
# convert bracket style from (...) to [...] inside the ln() wrapper.
$ws.Range("C2").Value  = "ln(GDP [dollars per capita])"
$ws.Range("C17").Value = "ln(GDP [dollars per capita])"

$ws.Range("C4").Value  = "ln(Tourism - Inbound [per capita])"

$ws.Range("C6").Value  = "ln(ProMed Mentions [per capita])"
$ws.Range("C13").Value = "ln(ProMed Mentions [per capita])"

$ws.Range("C8").Value  = "ln(Migrant Population [per capita])"

$ws.Range("C10").Value = "ln(AB Exports [dollars per capita])"

$ws.Range("C11").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C15").Value = "ln(Publication Bias Index [per capita])"

$ws.Range("C12").Value = "Livestock AB Consumption [kg per capita)"

$ws.Range("C16").Value = "ln(Population)"
